# CHC-identification.xlsx edit script
# Commit message: "rename first page of excel book"
#
# Applies the meaningful, content-level changes captured in the diff:
#   1. Rename the first worksheet "Identification" -> "CHC-Identification"
#   2. Move/restore the saved selection on that sheet to cell B15
#   3. Drop the workbook structure-lock protection (<workbookProtection> removed)

$wb = $excel.ActiveWorkbook

# 1) Rename the first worksheet.
$ws = $wb.Worksheets.Item(1)
$ws.Name = "CHC-Identification"

# 2) Restore the author's last selection on that sheet (was R11, now B15).
$ws.Activate()
$ws.Range("B15").Select()

# 3) Remove workbook (structure) protection/lock recorded in workbook.xml.
$wb.Unprotect()
